$wb = $excel.ActiveWorkbook

# --- Sheet 1: LoginCLU_Test_Cases_OK ---
# Clear out the two extra filled-in test case blocks (rows 4-9), leaving
# the first block (rows 1-3) untouched.
$ws1 = $wb.Worksheets.Item("LoginCLU_Test_Cases_OK")
$ws1.Range("A4:B9").ClearContents()
$ws1.Range("A4:XFD9").Select()

# --- Sheet 3: Consultas_Test_Cases_OK ---
# Remove the first ("CIF") test case block entirely (rows 1-3), which
# shifts the second ("Residencial Pospago") block up to become the new
# rows 1-3.
$ws3 = $wb.Worksheets.Item("Consultas_Test_Cases_OK")
$ws3.Rows("1:3").Delete()
$ws3.Range("A4:B9").ClearContents()
$ws3.Range("B13").Select()
